$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 229 (weekly price update), pushing the
# existing rows 229:249 down to 230:250. Excel's row-insert copies the
# formatting (incl. the date number-format on column D) from the row
# above, matching the rest of the table.
$ws.Rows.Item(229).EntireRow.Insert()

# Populate the newly inserted row with this week's observation. All of
# the descriptive/category columns repeat the same "Plátano / Sin
# especificar / Pintón / Ecuador" combination already used elsewhere in
# the sheet; only the date and the price columns are new.
$ws.Cells.Item(229, 1).Value = 1
$ws.Cells.Item(229, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(229, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(229, 4).Value = 44714
$ws.Cells.Item(229, 5).Value = 15
$ws.Cells.Item(229, 6).Value = "Fruta"
$ws.Cells.Item(229, 7).Value = 100108
$ws.Cells.Item(229, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(229, 9).Value = 100108006
$ws.Cells.Item(229, 10).Value = "Plátano"
$ws.Cells.Item(229, 11).Value = "Sin especificar"
$ws.Cells.Item(229, 12).Value = "Pintón"
$ws.Cells.Item(229, 13).Value = 120
$ws.Cells.Item(229, 14).Value = 14000
$ws.Cells.Item(229, 15).Value = 15000
$ws.Cells.Item(229, 16).Value = 14500
$ws.Cells.Item(229, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(229, 18).Value = "Ecuador"
$ws.Cells.Item(229, 19).Value = 725
$ws.Cells.Item(229, 20).Value = 20
